$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.823.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.41%  "

$ws.Range("D3").Value = "'1.657.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.14%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'311.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").Value = "'0.3636"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.88%  "

$ws.Range("D8").Value = "'47.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.67%  "

$ws.Range("D9").Value = "'0.3251"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.61%  "

$ws.Range("D10").Value = "'1.129"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.75%  "

$ws.Range("D11").Value = "'0.07051"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.52%  "

$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("D13").Value = "'6.040"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.26%  "

$ws.Range("D14").Value = "'19.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.00%  "

$ws.Range("D15").Value = "'1.654.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.03%  "

$ws.Range("D16").Value = "'6.574"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.50%  "

$ws.Range("D17").Value = "'0.00001045"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.72%  "

$ws.Range("D18").Value = "'0.06569"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.55%  "

$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").Value = "'78.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.94%  "

$ws.Range("D21").Value = "'5.887"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.29%  "

$ws.Range("D22").Value = "'15.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.69%  "

$ws.Range("D23").Value = "'12.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.23%  "

$ws.Range("D24").Value = "'24.840.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.24%  "

$ws.Range("D25").Value = "'2.445"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.49%  "

$ws.Range("D26").Value = "'2.442"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -9.11%  "

$ws.Range("D27").Value = "'147.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.72%  "

$ws.Range("D28").Value = "'18.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.38%  "

$ws.Range("D29").Value = "'1.838.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.00%  "

$ws.Range("D30").Value = "'1.197"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.89%  "

$ws.Range("D31").Value = "'125.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.03%  "

$ws.Range("D32").Value = "'4.085"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.58%  "

$ws.Range("D33").Value = "'5.747"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -11.29%  "

$ws.Range("D34").Value = "'0.08411"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.12%  "

$ws.Range("D35").Value = "'1.652"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.34%  "

$ws.Range("D36").Value = "'12.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.07%  "

$ws.Range("D37").Value = "'1.287"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.01%  "

$ws.Range("D38").Value = "'5.151"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.88%  "

$ws.Range("D39").Value = "'0.02248"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.74%  "

$ws.Range("D40").Value = "'0.06028"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.62%  "

$ws.Range("D41").Value = "'8.263"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.25%  "

$ws.Range("D42").Value = "'0.2059"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.35%  "

$ws.Range("D43").Value = "'1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").Value = "'0.5912"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.48%  "

$ws.Range("D45").Value = "'3.769"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.69%  "

$ws.Range("D46").Value = "'12.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.06%  "

$ws.Range("D47").Value = "'0.5600"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.04%  "

$ws.Range("D48").Value = "'124.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.53%  "

$ws.Range("D49").Value = "'1.938"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.85%  "

$ws.Range("D50").Value = "'0.06977"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.50%  "

$ws.Range("D51").Value = "'1.186"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.25%  "
